# "Modified Sign In Screen with DDD Framework!" - 2020/03/29
#
# Sheet1 ("Sign In Module") is reworked:
#   - the 3 old doLogin xpaths are replaced by the new loginid/password xpaths
#   - row 4 becomes "Sign In Submit Button" and a new row 5 carries the old
#     "Sign In Button" row down, with a new "button_signin" control id
#   - a small "Test Data" block is added at F1:G3 with a live mailto hyperlink
#     (styled the same way as the existing "Test Data" block on Contact Us)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sign In Module")
$contactUs = $wb.Worksheets.Item("Contact Us")

# --- Row 2: Sign In Email Address -> now targets the "loginid" control ---
$ws.Range("C2").Value2 = "loginid"
$ws.Range("D2").Value2 = "//input[@id='loginid']"

# --- Row 3: Sign In Password -> now targets the "password" control ---
$ws.Range("C3").Value2 = "password"
$ws.Range("D3").Value2 = "//input[@id='password']"

# --- Row 4: renamed to "Sign In Submit Button" with a new xpath ---
$ws.Range("A4").Value2 = "Sign In Submit Button"
$ws.Range("D4").Value2 = "//button[@class='btn btn-primary smtBtn ajax']"

# --- Row 5 (new): old "Sign In Button" row, now with a control id ---
$ws.Range("A5").Value2 = "Sign In Button"
$ws.Range("C5").Value2 = "button_signin"

# --- Test Data block (F1:G3), matching the look of the Contact Us block ---
$ws.Range("F1").Value2 = "Test Data"
$ws.Range("F2").Value2 = "email address"
$ws.Range("G2").Value2 = "stayseated05@gmail.com"
$ws.Range("F3").Value2 = "password"
$ws.Range("G3").Value2 = "1111111A"

$contactUs.Range("A10:B10").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:stayseated05@gmail.com") | Out-Null

$contactUs.Range("B11").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null

# --- Column widths (approximate autofit for the new/expanded columns) ---
$ws.Columns.Item(1).ColumnWidth = 19.5
$ws.Columns.Item(3).ColumnWidth = 12.67
$ws.Columns.Item(4).ColumnWidth = 43
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Columns.Item(7).ColumnWidth = 23.17
$ws.Columns.Item(8).ColumnWidth = 8.33

# --- Selection left on Sign In Module, but keep Contact Us as the active tab ---
$ws.Range("B14").Select() | Out-Null
$contactUs.Range("D7").Select() | Out-Null
$contactUs.Activate() | Out-Null

Write-Host "Sign In Module updated"
